$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row before row 3 (the existing "Peanut Candy" row),
# pushing it down to row 4 and carrying its formatting along.
$ws.Rows("3:3").Insert()

# --- New row 3: "Assassino Adryan desdobramento" ---
$ws.Range("A3").NumberFormat = "#,##0"
$ws.Range("A3").HorizontalAlignment = -4152
$ws.Range("B3").Value = "Assassino Adryan desdobramento"
$ws.Range("C3").Value = "noticias/Assassino-Adryan-desdobramento.html"
$ws.Range("D3").HorizontalAlignment = -4131
$ws.Range("D3").NumberFormat = "mm-dd-yy"
$ws.Range("D3").Value = 45743
$ws.Range("E3").Value = "images/Peanut-Candy-assassinato/images-2/suspeito-crime.png"

# --- Update the named range / table to cover the extra row ---
$wb.Names.Item("Notícias.accdb").RefersToR1C1 = "=Noticias!R1C1:R4C5"

Write-Host "done"
